$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column B (log_time) - shifts page_id..key_count right by one
$ws.Columns("B:B").Insert()

# 2. Set the new header
$ws.Range("B1").Value = "log_time"

# 3. Fill in the new log_time values for each row
$ws.Range("B2").Value = "2025-06-04 21:45:55"
$ws.Range("B3").Value = "2025-06-04 21:45:55"
$ws.Range("B4").Value = "2025-06-04 21:45:55"
$ws.Range("B5").Value = "2025-06-04 21:45:59"
$ws.Range("B6").Value = "2025-06-04 21:46:16"
$ws.Range("B7").Value = "2025-06-04 21:46:21"
$ws.Range("B8").Value = "2025-06-04 21:46:21"
$ws.Range("B9").Value = "2025-06-04 21:46:25"

# 4. Improve de-duplication logic: click_text is now also included as a key in
#    keys_combined / values_combined, bumping key_count by one for each row
#    that had a click_text value (rows with empty click_text stay unchanged).

# Row 2
$ws.Range("G2").Value = "channel, mainTitle, click_text, list_index, challengeName, challengeSeq, list_title, chal_index, activeParticipantCount, sticker"
$ws.Range("H2").Value = "Rround, 업로드, 3일차 이런 챌린지 어때요?, 0, 업로드, Optional(104), 업로드, 0, Optional(19), RECOMMEND"
$ws.Range("I2").Value = 10

# Row 3
$ws.Range("G3").Value = "channel, mainTitle, click_text, list_index, challengeName, challengeSeq, list_title, chal_index, activeParticipantCount, sticker"
$ws.Range("H3").Value = "Rround, 업로드, 이런 챌린지 어때요?, 0, 업로드, Optional(104), 업로드, 0, Optional(19), RECOMMEND"
$ws.Range("I3").Value = 10

# Row 4 (no click_text, unchanged besides the column shift already done by Insert)

# Row 5
$ws.Range("G5").Value = "channel, click_text, cta_text"
$ws.Range("H5").Value = "Rround, 클릭 텍스트, 인증하기"
$ws.Range("I5").Value = 3

# Row 6
$ws.Range("G6").Value = "click_text, goodsId, prd_name, channel, cta_text"
$ws.Range("H6").Value = "클릭 텍스트, 2655, 660, 정관장 홍삼대정 (홍삼대정 250g * 3병), [델리스푼] 브이핏 프리미엄 이너뷰티, Rround, 상품 선택 완료"
$ws.Range("I6").Value = 5

# Row 7 (no click_text, unchanged besides the column shift already done by Insert)

# Row 8 (no click_text, unchanged besides the column shift already done by Insert)

# Row 9
$ws.Range("G9").Value = "channel, click_text, popup_msg, cta_text, popup_title"
$ws.Range("H9").Value = "Rround, 클릭 텍스트, 다른 챌린지도 인증하고베스트 챌린저에 도전하세요 👏, 확인, 인증 완료!"
$ws.Range("I9").Value = 5
